$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H76").Value = 3558.2104
$ws_ALC.Range("J76").Value = 3550.3
$ws_ALC.Range("L76").Value = 3550.3
$ws_ALC.Range("N76").Value = -4180.3

$ws_ALC.Range("H79").Value = 3558.2104
$ws_ALC.Range("J79").Value = 3550.3
$ws_ALC.Range("L79").Value = 3550.3
$ws_ALC.Range("N79").Value = -5734.3

$ws_ALC.Range("H100").Value = 9259866
$ws_ALC.Range("I100").Value = 9615957
$ws_ALC.Range("J100").Value = 1500
$ws_ALC.Range("K100").Value = 9615957
$ws_ALC.Range("L100").Value = 1500
$ws_ALC.Range("M100").Value = -9615416
$ws_ALC.Range("N100").Value = -2582

$ws_ALC.Range("H132").Value = 15882578
$ws_ALC.Range("I132").Value = 23815868
$ws_ALC.Range("K132").Value = 71447604
$ws_ALC.Range("M132").Value = -71445074

$ws_ALC.Range("H137").Value = 1054.75
$ws_ALC.Range("I137").Value = 1028.9
$ws_ALC.Range("J137").Value = 1184
$ws_ALC.Range("K137").Value = 3086.7
$ws_ALC.Range("L137").Value = 3552
$ws_ALC.Range("M137").Value = -536.7000000000003
$ws_ALC.Range("N137").Value = -8652

$ws_ALC.Range("H141").Value = 783
$ws_ALC.Range("I141").Value = 783
$ws_ALC.Range("K141").Value = 2349
$ws_ALC.Range("M141").Value = 2831

$ws_ARM.Range("H61").Value = 2012.5
$ws_ARM.Range("I61").Value = 1350
$ws_ARM.Range("J61").Value = 4000
$ws_ARM.Range("K61").Value = 1350
$ws_ARM.Range("L61").Value = 4000
$ws_ARM.Range("M61").Value = -1138
$ws_ARM.Range("N61").Value = -4424

$ws_ARM.Range("H63").Value = 2272.375
$ws_ARM.Range("I63").Value = 2272.375
$ws_ARM.Range("K63").Value = 2272.375
$ws_ARM.Range("M63").Value = -1586.375

$ws_ARM.Range("H66").Value = 2272.375
$ws_ARM.Range("I66").Value = 2272.375
$ws_ARM.Range("K66").Value = 11361.875
$ws_ARM.Range("M66").Value = -7929.875

$ws_ARM.Range("H97").Value = 332.375
$ws_ARM.Range("I97").Value = 279.91666
$ws_ARM.Range("K97").Value = 279.91666
$ws_ARM.Range("M97").Value = 216.08334

$ws_ARM.Range("H110").Value = 1852
$ws_ARM.Range("I110").Value = 1310.1111
$ws_ARM.Range("J110").Value = 2339.7
$ws_ARM.Range("K110").Value = 1310.1111
$ws_ARM.Range("L110").Value = 2339.7
$ws_ARM.Range("M110").Value = 734.8888999999999
$ws_ARM.Range("N110").Value = -6429.7

$ws_ARM.Range("H122").Value = 1636.2307
$ws_ARM.Range("I122").Value = 1636.2307
$ws_ARM.Range("K122").Value = 4908.6921
$ws_ARM.Range("M122").Value = -2458.6921

$ws_ARM.Range("H124").Value = 16285.667
$ws_ARM.Range("J124").Value = 16285.667
$ws_ARM.Range("L124").Value = 16285.667
$ws_ARM.Range("N124").Value = -26105.667

$ws_ARM.Range("H125").Value = 44799.5
$ws_ARM.Range("J125").Value = 44799.5
$ws_ARM.Range("L125").Value = 44799.5
$ws_ARM.Range("N125").Value = -54639.5

$ws_ARM.Range("H132").Value = 2965.875
$ws_ARM.Range("I132").Value = 2880.077
$ws_ARM.Range("K132").Value = 8640.231
$ws_ARM.Range("M132").Value = -6110.231

$ws_ARM.Range("H136").Value = 2012.5
$ws_ARM.Range("I136").Value = 1350
$ws_ARM.Range("J136").Value = 4000
$ws_ARM.Range("K136").Value = 4050
$ws_ARM.Range("L136").Value = 12000
$ws_ARM.Range("M136").Value = -1500
$ws_ARM.Range("N136").Value = -17100

$ws_BSM.Range("H92").Value = 20998
$ws_BSM.Range("J92").Value = 20998
$ws_BSM.Range("L92").Value = 20998
$ws_BSM.Range("N92").Value = -25990

$ws_BSM.Range("H105").Value = 76925480
$ws_BSM.Range("I105").Value = 83335690
$ws_BSM.Range("K105").Value = 83335690
$ws_BSM.Range("M105").Value = -83333943

$ws_BSM.Range("H107").Value = 1841.3
$ws_BSM.Range("I107").Value = 1393.3077
$ws_BSM.Range("J107").Value = 2673.2856
$ws_BSM.Range("K107").Value = 1393.3077
$ws_BSM.Range("L107").Value = 2673.2856
$ws_BSM.Range("M107").Value = 526.6922999999999
$ws_BSM.Range("N107").Value = -6513.2856

$ws_BSM.Range("H134").Value = 7061.0527
$ws_BSM.Range("I134").Value = 1677.4
$ws_BSM.Range("K134").Value = 5032.200000000001
$ws_BSM.Range("M134").Value = -2497.200000000001

$ws_CRP.Range("H31").Value = 2228.923
$ws_CRP.Range("I31").Value = 1186.909
$ws_CRP.Range("J31").Value = 2993.0667
$ws_CRP.Range("K31").Value = 1186.909
$ws_CRP.Range("L31").Value = 2993.0667
$ws_CRP.Range("M31").Value = -891.9090000000001
$ws_CRP.Range("N31").Value = -3583.0667

$ws_CRP.Range("H34").Value = 2228.923
$ws_CRP.Range("I34").Value = 1186.909
$ws_CRP.Range("J34").Value = 2993.0667
$ws_CRP.Range("K34").Value = 1186.909
$ws_CRP.Range("L34").Value = 2993.0667
$ws_CRP.Range("M34").Value = -984.9090000000001
$ws_CRP.Range("N34").Value = -3397.0667

$ws_CRP.Range("H56").Value = 16924
$ws_CRP.Range("I56").Value = 16997
$ws_CRP.Range("J56").Value = 16899.666
$ws_CRP.Range("K56").Value = 16997
$ws_CRP.Range("L56").Value = 16899.666
$ws_CRP.Range("M56").Value = -16152
$ws_CRP.Range("N56").Value = -18589.666

$ws_CRP.Range("H58").Value = 975.0625
$ws_CRP.Range("I58").Value = 973.4
$ws_CRP.Range("J58").Value = 1000
$ws_CRP.Range("K58").Value = 973.4
$ws_CRP.Range("L58").Value = 1000
$ws_CRP.Range("M58").Value = -770.4
$ws_CRP.Range("N58").Value = -1406

$ws_CRP.Range("H111").Value = 30234
$ws_CRP.Range("J111").Value = 30234
$ws_CRP.Range("L111").Value = 30234
$ws_CRP.Range("N111").Value = -38414

$ws_CRP.Range("H132").Value = 2450.8462
$ws_CRP.Range("I132").Value = 2057.75
$ws_CRP.Range("J132").Value = 3079.8
$ws_CRP.Range("K132").Value = 6173.25
$ws_CRP.Range("L132").Value = 9239.400000000001
$ws_CRP.Range("M132").Value = -3643.25
$ws_CRP.Range("N132").Value = -14299.4

$ws_CRP.Range("H134").Value = 18519570
$ws_CRP.Range("I134").Value = 22223250
$ws_CRP.Range("J134").Value = 1171.3334
$ws_CRP.Range("K134").Value = 66669750
$ws_CRP.Range("L134").Value = 3514.0002
$ws_CRP.Range("M134").Value = -66667215
$ws_CRP.Range("N134").Value = -8584.0002

$ws_CRP.Range("H136").Value = 975.0625
$ws_CRP.Range("I136").Value = 973.4
$ws_CRP.Range("J136").Value = 1000
$ws_CRP.Range("K136").Value = 2920.2
$ws_CRP.Range("L136").Value = 3000
$ws_CRP.Range("M136").Value = -370.1999999999998
$ws_CRP.Range("N136").Value = -8100

$ws_CUL.Range("H33").Value = 405.66666
$ws_CUL.Range("I33").Value = 100
$ws_CUL.Range("J33").Value = 466.8
$ws_CUL.Range("K33").Value = 600
$ws_CUL.Range("L33").Value = 2800.8
$ws_CUL.Range("M33").Value = -317
$ws_CUL.Range("N33").Value = -3366.8

$ws_CUL.Range("H68").Value = 960
$ws_CUL.Range("J68").Value = 1067.6666
$ws_CUL.Range("L68").Value = 3202.9998
$ws_CUL.Range("N68").Value = -4824.9998

$ws_CUL.Range("H71").Value = 960
$ws_CUL.Range("J71").Value = 1067.6666
$ws_CUL.Range("L71").Value = 9608.999400000001
$ws_CUL.Range("N71").Value = -17720.9994

$ws_CUL.Range("H107").Value = 790.8
$ws_CUL.Range("I107").Value = 0
$ws_CUL.Range("J107").Value = 790.8
$ws_CUL.Range("K107").Value = 0
$ws_CUL.Range("L107").Value = 2372.4
$ws_CUL.Range("N107").Value = -6212.4
$ws_CUL.Range("M107").ClearContents()

$ws_CUL.Range("H140").Value = 1989.0834
$ws_CUL.Range("J140").Value = 2890
$ws_CUL.Range("L140").Value = 8670
$ws_CUL.Range("N140").Value = -19030

$ws_GSM.Range("H80").Value = 2507.1428
$ws_GSM.Range("I80").Value = 1583.3334
$ws_GSM.Range("J80").Value = 3200
$ws_GSM.Range("K80").Value = 1583.3334
$ws_GSM.Range("L80").Value = 3200
$ws_GSM.Range("M80").Value = -585.3334
$ws_GSM.Range("N80").Value = -5196

$ws_GSM.Range("H83").Value = 2507.1428
$ws_GSM.Range("I83").Value = 1583.3334
$ws_GSM.Range("J83").Value = 3200
$ws_GSM.Range("K83").Value = 7916.666999999999
$ws_GSM.Range("L83").Value = 16000
$ws_GSM.Range("M83").Value = -2924.666999999999
$ws_GSM.Range("N83").Value = -25984

$ws_GSM.Range("H107").Value = 816.8
$ws_GSM.Range("I107").Value = 789.1111
$ws_GSM.Range("J107").Value = 858.3333
$ws_GSM.Range("K107").Value = 789.1111
$ws_GSM.Range("L107").Value = 858.3333
$ws_GSM.Range("M107").Value = 1130.8889
$ws_GSM.Range("N107").Value = -4698.3333

$ws_GSM.Range("H109").Value = 20000
$ws_GSM.Range("J109").Value = 20000
$ws_GSM.Range("L109").Value = 20000
$ws_GSM.Range("N109").Value = -22080

$ws_GSM.Range("H113").Value = 2995.3333
$ws_GSM.Range("I113").Value = 1316.1428
$ws_GSM.Range("K113").Value = 1316.1428
$ws_GSM.Range("M113").Value = 853.8571999999999

$ws_GSM.Range("H122").Value = 215957.14
$ws_GSM.Range("I122").Value = 1066.6666
$ws_GSM.Range("J122").Value = 377125
$ws_GSM.Range("K122").Value = 3199.9998
$ws_GSM.Range("L122").Value = 1131375
$ws_GSM.Range("M122").Value = -749.9998000000001
$ws_GSM.Range("N122").Value = -1136275

$ws_GSM.Range("H132").Value = 2407.3333
$ws_GSM.Range("I132").Value = 1736.25
$ws_GSM.Range("J132").Value = 3749.5
$ws_GSM.Range("K132").Value = 5208.75
$ws_GSM.Range("L132").Value = 11248.5
$ws_GSM.Range("M132").Value = -2678.75
$ws_GSM.Range("N132").Value = -16308.5

$ws_LTW.Range("H93").Value = 0
$ws_LTW.Range("I93").Value = 0
$ws_LTW.Range("J93").Value = 0
$ws_LTW.Range("K93").Value = 0
$ws_LTW.Range("L93").Value = 0
$ws_LTW.Range("M93").ClearContents()
$ws_LTW.Range("N93").ClearContents()

$ws_LTW.Range("H122").Value = 27781068
$ws_LTW.Range("I122").Value = 83337300
$ws_LTW.Range("K122").Value = 250011900
$ws_LTW.Range("M122").Value = -250009450

$ws_LTW.Range("H132").Value = 78938.92
$ws_LTW.Range("I132").Value = 1212.875
$ws_LTW.Range("J132").Value = 203300.6
$ws_LTW.Range("K132").Value = 3638.625
$ws_LTW.Range("L132").Value = 609901.8
$ws_LTW.Range("M132").Value = -1108.625
$ws_LTW.Range("N132").Value = -614961.8

$ws_WVR.Range("H96").Value = 2222.7144
$ws_WVR.Range("I96").Value = 1976
$ws_WVR.Range("K96").Value = 1976
$ws_WVR.Range("M96").Value = -603

$ws_WVR.Range("H113").Value = 942.1429000000001
$ws_WVR.Range("I113").Value = 298.75
$ws_WVR.Range("K113").Value = 896.25
$ws_WVR.Range("M113").Value = 1273.75

$ws_WVR.Range("H132").Value = 3614.5217
$ws_WVR.Range("I132").Value = 3131.5293
$ws_WVR.Range("K132").Value = 9394.5879
$ws_WVR.Range("M132").Value = -6864.5879
